# Auto-generated Excel COM-interop edit script
# Applies cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need to be
# forced to Text format first, otherwise Excel auto-converts the string
# into a numeric value (losing the original text formatting/precision).
$textForceCells = @(
  "D5",
  "D6",
  "D10",
  "D12",
  "D13",
  "D14",
  "D18",
  "D21",
  "D22",
  "D23",
  "D24",
  "D25",
  "D26",
  "D27",
  "D29",
  "D30",
  "D31",
  "D33",
  "D34",
  "D36",
  "D37",
  "D39",
  "D42",
  "D44",
  "D45",
  "D47",
  "D48",
  "D49",
  "D50"
)
foreach ($addr in $textForceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "68.233.33"
$ws.Range("D3").Value = "3.682.19"
$ws.Range("E3").Value = "  +6.08%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "421.39"
$ws.Range("E5").Value = "  +1.66%  "
$ws.Range("D6").Value = "130.29"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("D8").Value = "3.671.68"
$ws.Range("E8").Value = "  +5.92%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.769"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("E11").Value = "  +27.28%  "
$ws.Range("D12").Value = "0.0000446"
$ws.Range("E12").Value = "  +96.89%  "
$ws.Range("D13").Value = "42.04"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "9.81"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").Value = "4.244.74"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("E16").Value = "  +0.36%  "
$ws.Range("D17").Value = "3.677.94"
$ws.Range("E17").Value = "  +6.32%  "
$ws.Range("D18").Value = "20.12"
$ws.Range("E18").Value = "  -0.40%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "68.159.49"
$ws.Range("E20").Value = "  +7.44%  "
$ws.Range("D21").Value = "12.53"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").Value = "460.45"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "89.16"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "13.59"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").Value = "  -7.08%  "
$ws.Range("D26").Value = "36.75"
$ws.Range("E26").Value = "  +10.12%  "
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "4.97"
$ws.Range("E29").Value = "  +4.51%  "
$ws.Range("D30").Value = "2.80"
$ws.Range("E30").Value = "  +4.95%  "
$ws.Range("D31").Value = "12.25"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("E32").Value = "  +6.01%  "
$ws.Range("D33").Value = "7.12"
$ws.Range("E33").Value = "  -5.05%  "
$ws.Range("D34").Value = "40.27"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  -7.12%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "56.14"
$ws.Range("E37").Value = "  -2.60%  "
$ws.Range("D38").Value = "0.0₃0790"
$ws.Range("E38").Value = "  +22.03%  "
$ws.Range("D39").Value = "0.0491"
$ws.Range("E39").Value = "  +1.36%  "
$ws.Range("E40").Value = "  +10.38%  "
$ws.Range("D42").Value = "149.20"
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  -5.29%  "
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("E46").Value = "  +13.54%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "4.26"
$ws.Range("E47").Value = "  -5.40%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.168"
$ws.Range("E48").Value = "  +21.35%  "
$ws.Range("D49").Value = "0.300"
$ws.Range("E49").Value = "  -4.24%  "
$ws.Range("D50").Value = "1.96"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("E51").Value = "  +13.16%  "
